$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "19-03-2025" (F:G) and "20-03-2025" (H:I) columns.
# The remaining D:E columns (previously "18-03-2025") are repurposed below
# to hold the new "20-03-2025" header/data.
$ws.Range("F1:I6").EntireColumn.Delete()

# Update the D/E headers to reflect the new date.
$ws.Range("D1").Value = "20-03-2025 Status"
$ws.Range("E1").Value = "20-03-2025 Time"

# Refresh the attendance totals and the 20-03-2025 check-in times.
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = "13:48:06"

$ws.Range("C3").Value = 1
$ws.Range("E3").Value = "13:48:15"

$ws.Range("C4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("C6").Value = ""
